$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.842.31"
$ws.Range("D3").Value = "1.763.28"
$ws.Range("E3").Value = "  +1.58%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.27"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4460"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3541"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07401"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.90"
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.098"
$ws.Range("E11").Value = "  +2.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.86"
$ws.Range("E13").Value = "  +2.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.015"
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.236"
$ws.Range("E15").Value = "  +2.88%  "
$ws.Range("D16").Value = "1.761.97"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.85"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001060"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06430"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("E21").Value = "  +3.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.758"
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("D23").Value = "27.876.87"
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.23"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.107"
$ws.Range("E25").Value = "  +1.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.61"
$ws.Range("E26").Value = "  -0.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.33"
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("D28").Value = "1.965.75"
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.141"
$ws.Range("E29").Value = "  +5.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.08"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.098"
$ws.Range("E31").Value = "  +5.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09192"
$ws.Range("E32").Value = "  +1.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.641"
$ws.Range("E33").Value = "  +5.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.692"
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.82"
$ws.Range("E35").Value = "  +2.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06179"
$ws.Range("E36").Value = "  +4.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02275"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2098"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6306"
$ws.Range("E39").Value = "  +1.47%  "
$ws.Range("E40").Value = "  +1.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.181"
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.392"
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.829"
$ws.Range("E43").Value = "  +2.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.24"
$ws.Range("E44").Value = "  +1.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.736"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5840"
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.27"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.950"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06876"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.134"
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.73"
$ws.Range("E51").Value = "  +2.50%  "
